# C5-PowerPoint.pptx edit — Wed, Apr 08, 2020
#
# 1) Slide 6's table switches to a different built-in table style
#    ({AF8019D1-8A22-4B26-AF64-7CAA1E32C1E7} -> {3DA207B8-406E-4D1D-A55F-C7A5E3F7F15F}).
# 2) The deck's theme colour scheme (currently the "Integral" palette) is swapped
#    for the stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Retarget the table style on slide 6 -------------------------------
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table
$tbl.ApplyStyle("{3DA207B8-406E-4D1D-A55F-C7A5E3F7F15F}")

# --- 2. Swap the theme colour scheme from "Integral" to "Office" ----------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Each value is the 0x00BBGGRR long used by ThemeColorScheme.Colors(i).RGB
$tcs.Colors(1).RGB  = 0x000000    # dk1      000000
$tcs.Colors(2).RGB  = 0xFFFFFF    # lt1      FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444    # dk2      44546A
$tcs.Colors(4).RGB  = 0xE6E6E7    # lt2      E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B    # accent1  5B9BD5
$tcs.Colors(6).RGB  = 0x317DED    # accent2  ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5    # accent3  A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF    # accent4  FFC000
$tcs.Colors(9).RGB  = 0xC47244    # accent5  4472C4
$tcs.Colors(10).RGB = 0x47AD70    # accent6  70AD47
$tcs.Colors(11).RGB = 0xC16305    # hyperlink         0563C1
$tcs.Colors(12).RGB = 0x724F95    # followed hyperlink 954F72
